$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.349.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.03%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.766.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.07%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4299"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07079"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8472"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.25"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.732.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.81%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.432"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06788"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008653"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.354.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.021"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.992.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.68"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.863"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.070"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.69%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.715"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.23%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7317"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.336"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.33%  "

# Row 34
$ws.Range("E34").Value = "  -0.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.076"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05132"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01888"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4915"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1609"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "

# Row 42
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.515"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.47%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.064"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "

# Row 47
$ws.Range("E47").Value = "  -3.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4495"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06192"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.578"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.730"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
